$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update statistic values for several existing countries (COVID data refresh) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1187233
$ws.Range("C4").Value = 26459
$ws.Range("D4").Value = 178263
$ws.Range("E4").Value = 940404
$ws.Range("G4").Value = 1122
$ws.Range("H4").Value = 68566

# Row 60: Kazajistan
$ws.Range("B60").Value = 3920
$ws.Range("C60").Value = 63
$ws.Range("E60").Value = 2809

# Row 68: Nigeria
$ws.Range("B68").Value = 2558
$ws.Range("C68").Value = 170
$ws.Range("D68").Value = 400
$ws.Range("E68").Value = 2071
$ws.Range("F68").Value = 4
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 87

# Row 99: Principado de Andorra
$ws.Range("B99").Value = 748
$ws.Range("C99").Value = 1
$ws.Range("D99").Value = 493
$ws.Range("E99").Value = 210
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 45

# Row 108: Uruguay
$ws.Range("B108").Value = 655
$ws.Range("C108").Value = 3
$ws.Range("D108").Value = 442
$ws.Range("E108").Value = 196

# Row 153: Bermudas
$ws.Range("B153").Value = 115
$ws.Range("C153").Value = 1
$ws.Range("E153").Value = 57

# --- Re-position "Santa Lucia": it moves from its old slot (after Gambia) to a new
# slot right after "Fiyi" (before "Belice"), with refreshed statistics that now
# rank it higher in the table (tied with Fiyi/Belice/Nueva Caledonia at 18 total cases).
# Rows 188-192 shift down by one to make room, and the old "Santa Lucia" row is dropped.

$ws.Range("A188").Value = "Santa Lucia"
$ws.Range("B188").Value = 18
$ws.Range("C188").Value = 1
$ws.Range("D188").Value = 15
$ws.Range("E188").Value = 3
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 0

$ws.Range("A189").Value = "Belice"
$ws.Range("B189").Value = 18
$ws.Range("C189").Value = 0
$ws.Range("D189").Value = 13
$ws.Range("E189").Value = 3
$ws.Range("F189").Value = 1
$ws.Range("G189").Value = 0
$ws.Range("H189").Value = 2

$ws.Range("A190").Value = "Nueva Caledonia"
$ws.Range("B190").Value = 18
$ws.Range("C190").Value = 0
$ws.Range("D190").Value = 17
$ws.Range("E190").Value = 1
$ws.Range("F190").Value = 1
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 0

$ws.Range("A191").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("B191").Value = 17
$ws.Range("C191").Value = 0
$ws.Range("D191").Value = 0
$ws.Range("E191").Value = 17
$ws.Range("F191").Value = 0
$ws.Range("G191").Value = 0
$ws.Range("H191").Value = 0

$ws.Range("A192").Value = "Gambia"
$ws.Range("B192").Value = 17
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 9
$ws.Range("E192").Value = 7
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 1

# Row 193 (Santo Tome y Principe) stays unchanged.
